$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) values
# to match the Oct 20 2023 13:58:55 UTC data pull.
# Leading apostrophe on column D assignments forces text storage so values
# like "213.11" or "1.602.53" are kept verbatim instead of being parsed as numbers.
$ws.Range("D2").Value = "'29.556.16"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").Value = "'1.602.53"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'213.11"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("E6").Value = "  +7.13%  "
$ws.Range("D8").Value = "'26.87"
$ws.Range("E8").Value = "  +10.89%  "
$ws.Range("D9").Value = "'43.50"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").Value = "'0.250"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "'0.0914"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").Value = "'1.831.26"
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("D14").Value = "'1.599.31"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("D15").Value = "'29.548.06"
$ws.Range("E15").Value = "  +3.78%  "
$ws.Range("D16").Value = "'0.536"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D18").Value = "'63.55"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("D19").Value = "'240.51"
$ws.Range("E19").Value = "  +4.91%  "
$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "'154.46"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("E28").Value = "  +5.03%  "
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").Value = "'1.427.91"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").Value = "'1.03"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "'2.81"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.0166"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").Value = "'0.535"
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'53.64"
$ws.Range("E43").Value = "  +23.00%  "
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "'0.996"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "'0.0471"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "'65.57"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "'1.741.95"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("D50").Value = "'86.53"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("E51").Value = "  -3.57%  "
